$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add new row 6 (Andres) first, so the shared strings table picks
# up the new unique values in the same order the original author
# introduced them.
# ---------------------------------------------------------------
$ws.Range("A6").Value = "CC"
$ws.Range("B6").Value = 1107834925

# C6 - nombre_completo, same text style as C2:C5 (s=1)
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Andrés"

# D6 / E6 - correo_electronico / correo_institucional -> hyperlinked email
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "andres@gmail.com"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "andres@gmail.com"

$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:andres@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:andres@gmail.com")

$ws.Range("F6").Value = 320

# G6 - fecha_nacimiento, same date style as G2:G5 (s=2)
$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value2 = 38054

$ws.Range("H6").Value = 3207282500

# I6 - promedio_acumulado, stored as text
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "4.2"
$ws.Range("A6").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$ws.Range("J6").Value = 30
$ws.Range("K6").Value = "M"
$ws.Range("L6").Value = "A00123300"

# ---------------------------------------------------------------
# Update existing rows
# ---------------------------------------------------------------

# Row 3: tipo_documento TI -> CC
$ws.Range("A3").Value = "CC"

# Column F (puntaje_icfes) updates
$ws.Range("F2").Value = 409
$ws.Range("F3").Value = 123
$ws.Range("F4").Value = 390
$ws.Range("F5").Value = 123

# Column I (promedio_acumulado) updates - keep stored as text
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "4.1"
$ws.Range("A2").Copy()
$ws.Range("I4").PasteSpecial(-4122)

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "3.9"
$ws.Range("A2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet view: drop the frozen/scrolled topLeftCell and move the
# active selection to A4
# ---------------------------------------------------------------
$ws.Range("A4").Select()

Write-Host "done"
